# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 dashboard snapshot. This refresh:
#   1. Bumps the "last updated" timestamp in A1 (09:47 -> 11:04).
#   2. Pulls fresh numbers for several countries, which also shuffles a
#      handful of countries up/down by one row because the sheet is kept
#      sorted by total cases (column B) descending. Each affected row's
#      country name (column A) and its B/C/D/E/G/H stats are rewritten to
#      match the new snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refreshed-at timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 11:04"

# --- Banglades / Alemania / Indonesia / Filipinas block (rows 20-23) --
# Indonesia overtakes Alemania.
$ws.Range("A21").Value = "Indonesia"
$ws.Range("B21").Value = 357762
$ws.Range("C21").Value = 4301
$ws.Range("D21").Value = 281592
$ws.Range("E21").Value = 63739
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 12431

$ws.Range("A22").Value = "Alemania"
$ws.Range("B22").Value = 356792
$ws.Range("D22").Value = 287600
$ws.Range("E22").Value = 59356
$ws.Range("H22").Value = 9836

$ws.Range("B23").Value = 354338
$ws.Range("C23").Value = 2673
$ws.Range("D23").Value = 295312
$ws.Range("E23").Value = 52423
$ws.Range("G23").Value = 73
$ws.Range("H23").Value = 6603

# --- Rumania / Marruecos / Chequia / Polonia block (rows 32-35) -------
# Polonia overtakes Marruecos and Chequia.
$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 167230
$ws.Range("C33").Value = 9622
$ws.Range("D33").Value = 90162
$ws.Range("E33").Value = 73544
$ws.Range("G33").Value = 84
$ws.Range("H33").Value = 3524

$ws.Range("A34").Value = "Marruecos"
$ws.Range("B34").Value = 167148
$ws.Range("D34").Value = 138989
$ws.Range("E34").Value = 25341
$ws.Range("H34").Value = 2818

$ws.Range("A35").Value = "Chequia"
$ws.Range("B35").Value = 160112
$ws.Range("D35").Value = 66093
$ws.Range("E35").Value = 92736
$ws.Range("H35").Value = 1283

# --- Singapur (row 65) --------------------------------------------------
$ws.Range("B65").Value = 57904
$ws.Range("C65").Value = 3
$ws.Range("E65").Value = 92

# --- Bulgaria / Australia / Eslovaquia block (rows 84-86) -------------
# Eslovaquia overtakes Australia.
$ws.Range("A85").Value = "Eslovaquia"
$ws.Range("B85").Value = 28268
$ws.Range("C85").Value = 1968
$ws.Range("D85").Value = 7297
$ws.Range("E85").Value = 20889
$ws.Range("G85").Value = 11
$ws.Range("H85").Value = 82

$ws.Range("A86").Value = "Australia"
$ws.Range("B86").Value = 27383
$ws.Range("C86").Value = 12
$ws.Range("D86").Value = 25098
$ws.Range("E86").Value = 1381
$ws.Range("H86").Value = 904

# --- Corea del Sur / Grecia / Croacia block (rows 87-89) ---------------
# Croacia overtakes Grecia.
$ws.Range("A88").Value = "Croacia"
$ws.Range("B88").Value = 24761
$ws.Range("C88").Value = 1096
$ws.Range("D88").Value = 19562
$ws.Range("E88").Value = 4844
$ws.Range("G88").Value = 10
$ws.Range("H88").Value = 355

$ws.Range("A89").Value = "Grecia"
$ws.Range("B89").Value = 24450
$ws.Range("D89").Value = 9989
$ws.Range("E89").Value = 13971
$ws.Range("H89").Value = 490

# --- Finlandia (row 102) ------------------------------------------------
$ws.Range("B102").Value = 13293
$ws.Range("C102").Value = 160
$ws.Range("E102").Value = 3842

# --- Consejo Danes para los Refugiados (row 107) ------------------------
$ws.Range("B107").Value = 11000
$ws.Range("C107").Value = 1
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 302

# --- Cabo Verde / Angola / Guadalupe / Lituania block (rows 118-121) ---
# Lituania overtakes Angola and Guadalupe.
$ws.Range("A119").Value = "Lituania"
$ws.Range("B119").Value = 7269
$ws.Range("C119").Value = 228
$ws.Range("D119").Value = 3097
$ws.Range("E119").Value = 4059
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 113

$ws.Range("A120").Value = "Angola"
$ws.Range("B120").Value = 7222
$ws.Range("D120").Value = 3012
$ws.Range("E120").Value = 3976
$ws.Range("H120").Value = 234

$ws.Range("A121").Value = "Guadalupe"
$ws.Range("B121").Value = 7122
$ws.Range("D121").Value = 2199
$ws.Range("E121").Value = 4827
$ws.Range("H121").Value = 96
